$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for be8735b1 row now reports
# the actual handoff datetime instead of the "in sync" status text.
$wsOverview.Range("G3").Value = "2016-08-31 10:52:48"

# zh-cn sheet: regenerate handoff/handback datetimes for the 58f57eca row.
$wsZhCn.Range("H2").Value = "2016-08-31 10:53:48"
$wsZhCn.Range("K2").Value = "2016-08-31 10:54:14"

# de-de sheet: regenerate handback datetime for the 58f57eca row, and the
# handoff datetime for the be8735b1 row.
$wsDeDe.Range("K2").Value = "2016-08-31 10:54:20"
$wsDeDe.Range("H3").Value = "2016-08-31 10:52:48"
